$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the CPMV (6-parameter motion vector) candidate values in row 2,
# which previously held the mismatched set of numbers.
$ws.Range("L2").Value = 2
$ws.Range("M2").Value = -16
$ws.Range("N2").Value = -80
$ws.Range("O2").Value = 80
$ws.Range("P2").Value = -240
$ws.Range("Q2").Value = 16
$ws.Range("R2").Value = -32

# Row 4's M:R values were the correct CPMV parameters that have now been
# consolidated into row 2 above, so clear them from row 4.
$ws.Range("M4:R4").ClearContents()

# Remove the now-unused comp_offs_x / comp_offs_y columns (U and V),
# including their header cells in row 1.
$ws.Range("U1:V3").ClearContents()

# Update the view state to match the saved selection.
$ws.Range("S6").Select()
